$d = $word.ActiveDocument

# --- 1. Insert three new bullet paragraphs before the
#        "Developed and deployed custom analytical tools..." bullet,
#        within the Siege Analytics / "Advanced Data Analysis and
#        Statistical Modeling" section. ---

$bullet = [char]0x2022

$findRange = $d.Content
$findRange.Find.Execute(
    "Developed and deployed custom analytical tools and algorithms using Python",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$targetPara = $findRange.Paragraphs(1)
$targetStart = $targetPara.Range.Start
$insertPoint = $d.Range($targetStart, $targetStart)

$newText1 = $bullet + " Developed meta-analytical techniques that identified systematic data quality issues across 20+ years of voter registration data`r"
$newText2 = $bullet + " Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters`r"
$newText3 = $bullet + " Created fraud detection systems analyzing 5+ terabyte datasets, uncovering demographic miscoding patterns across 2,000+ precincts`r"

$insertPoint.InsertBefore($newText1 + $newText2 + $newText3)

# --- 2. Remove the old "Created fraud detection systems for campaign
#        finance data analysis across multi-terabyte datasets" bullet
#        (it followed the "Developed meta-analytical techniques to
#        resolve ambiguous dimensions..." bullet). ---

$removeRange = $d.Content
$removeRange.Find.Execute(
    "Created fraud detection systems for campaign finance data analysis across multi-terabyte datasets",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$removePara = $removeRange.Paragraphs(1)
$removePara.Range.Delete()
